# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Rules!B11 ("R40") is replaced by the text "1" while keeping its
# existing cell style (s="23") and its "string" cell type (t="s").
#
# A plain `Range.Value = "1"` assignment would be auto-coerced by Excel
# into the *number* 1 (and picking up a different/new style in the
# process), which is not what we want here: the target cell must keep
# storing a shared-string "1", not a numeric 1. To reproduce a literal
# text value exactly like Excel does when a formula result is pasted as
# a value, we enter a text formula and then collapse it down to a plain
# value with Copy / PasteSpecial(xlPasteValues). This keeps the original
# formatting (xf/style) of the cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")
$target.Formula = "=""1"""
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
